$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price/volume snapshot (and a ranking swap between
# FraxShare and TrustWalletToken in rows 40-41). Each entry is the final
# text that belongs in the given cell. A leading "'" forces Excel to store
# a numeric-looking price (e.g. "4.419") as literal text instead of
# re-parsing it into a Double, matching the sheet's existing convention of
# keeping every Price/Volume column value as a string.
$updates = @(
    @{ Cell = 'D2'; Value = '26.097.70' },
    @{ Cell = 'E2'; Value = '  -0.28%  ' },
    @{ Cell = 'D3'; Value = '1.646.92' },
    @{ Cell = 'E3'; Value = '  -1.37%  ' },
    @{ Cell = 'E4'; Value = '  -0.15%  ' },
    @{ Cell = 'D5'; Value = '''216.30' },
    @{ Cell = 'E5'; Value = '  +2.72%  ' },
    @{ Cell = 'D6'; Value = '''0.5211' },
    @{ Cell = 'E6'; Value = '  +0.03%  ' },
    @{ Cell = 'E7'; Value = '  -0.17%  ' },
    @{ Cell = 'E8'; Value = '  -0.40%  ' },
    @{ Cell = 'D9'; Value = '''0.06362' },
    @{ Cell = 'E9'; Value = '  +0.64%  ' },
    @{ Cell = 'E10'; Value = '  -1.67%  ' },
    @{ Cell = 'D11'; Value = '''0.07671' },
    @{ Cell = 'E11'; Value = '  +1.59%  ' },
    @{ Cell = 'D12'; Value = '1.648.12' },
    @{ Cell = 'E12'; Value = '  -1.41%  ' },
    @{ Cell = 'D13'; Value = '''4.419' },
    @{ Cell = 'E13'; Value = '  -0.52%  ' },
    @{ Cell = 'D14'; Value = '1.869.30' },
    @{ Cell = 'D15'; Value = '''0.5547' },
    @{ Cell = 'E15'; Value = '  +0.97%  ' },
    @{ Cell = 'D16'; Value = '0.0₅8302' },
    @{ Cell = 'E16'; Value = '  +3.46%  ' },
    @{ Cell = 'D17'; Value = '''65.01' },
    @{ Cell = 'E17'; Value = '  -2.06%  ' },
    @{ Cell = 'D18'; Value = '26.104.30' },
    @{ Cell = 'E18'; Value = '  -0.29%  ' },
    @{ Cell = 'E19'; Value = '  -0.09%  ' },
    @{ Cell = 'D20'; Value = '''4.732' },
    @{ Cell = 'E20'; Value = '  -0.47%  ' },
    @{ Cell = 'D21'; Value = '''188.40' },
    @{ Cell = 'E21'; Value = '  +0.67%  ' },
    @{ Cell = 'D22'; Value = '''10.21' },
    @{ Cell = 'D23'; Value = '''6.228' },
    @{ Cell = 'E23'; Value = '  +0.26%  ' },
    @{ Cell = 'E24'; Value = '  -0.20%  ' },
    @{ Cell = 'D25'; Value = '''146.19' },
    @{ Cell = 'E25'; Value = '  -2.46%  ' },
    @{ Cell = 'D26'; Value = '''0.1218' },
    @{ Cell = 'E26'; Value = '  -1.89%  ' },
    @{ Cell = 'D27'; Value = '''7.432' },
    @{ Cell = 'E27'; Value = '  -0.78%  ' },
    @{ Cell = 'E28'; Value = '  +0.11%  ' },
    @{ Cell = 'D29'; Value = '''1.385' },
    @{ Cell = 'E29'; Value = '  +2.35%  ' },
    @{ Cell = 'D30'; Value = '''0.05978' },
    @{ Cell = 'E30'; Value = '  -5.70%  ' },
    @{ Cell = 'D31'; Value = '''1.271' },
    @{ Cell = 'E31'; Value = '  -0.92%  ' },
    @{ Cell = 'D32'; Value = '''3.408' },
    @{ Cell = 'E32'; Value = '  -3.21%  ' },
    @{ Cell = 'D33'; Value = '''3.397' },
    @{ Cell = 'E33'; Value = '  -0.50%  ' },
    @{ Cell = 'D34'; Value = '''1.662' },
    @{ Cell = 'E34'; Value = '  +1.10%  ' },
    @{ Cell = 'D35'; Value = '''0.9935' },
    @{ Cell = 'E35'; Value = '  -1.13%  ' },
    @{ Cell = 'E36'; Value = '  -0.47%  ' },
    @{ Cell = 'D37'; Value = '''2.754' },
    @{ Cell = 'E37'; Value = '  -0.02%  ' },
    @{ Cell = 'D38'; Value = '''0.5649' },
    @{ Cell = 'E38'; Value = '  -6.65%  ' },
    @{ Cell = 'D39'; Value = '''0.01618' },
    @{ Cell = 'E39'; Value = '  +0.22%  ' },
    @{ Cell = 'B40'; Value = 'TrustWalletToken' },
    @{ Cell = 'C40'; Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt' },
    @{ Cell = 'D40'; Value = '''0.8612' },
    @{ Cell = 'E40'; Value = '  -0.48%  ' },
    @{ Cell = 'B41'; Value = 'FraxShare' },
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs' },
    @{ Cell = 'D41'; Value = '''5.848' },
    @{ Cell = 'E41'; Value = '  -4.58%  ' },
    @{ Cell = 'D42'; Value = '''1.001' },
    @{ Cell = 'E42'; Value = '  -0.31%  ' },
    @{ Cell = 'D43'; Value = '1.029.16' },
    @{ Cell = 'E43'; Value = '  -7.49%  ' },
    @{ Cell = 'D44'; Value = '''99.99' },
    @{ Cell = 'E44'; Value = '  -0.39%  ' },
    @{ Cell = 'D45'; Value = '1.796.24' },
    @{ Cell = 'E45'; Value = '  -1.46%  ' },
    @{ Cell = 'D46'; Value = '0.0₈110' },
    @{ Cell = 'E46'; Value = '  +1.78%  ' },
    @{ Cell = 'D47'; Value = '''55.93' },
    @{ Cell = 'E47'; Value = '  +0.56%  ' },
    @{ Cell = 'D48'; Value = '''0.9998' },
    @{ Cell = 'E48'; Value = '  -0.06%  ' },
    @{ Cell = 'D49'; Value = '''8.055' },
    @{ Cell = 'E49'; Value = '  -0.39%  ' },
    @{ Cell = 'D50'; Value = '''0.05170' },
    @{ Cell = 'E50'; Value = '  -1.27%  ' },
    @{ Cell = 'E51'; Value = '  -0.47%  ' }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
